$d = $word.ActiveDocument

# Locate the paragraph that ends with "... 68 ... fixed by removing it"
# (the 2nd / last occurrence of the "fixed by removing it" paragraph in the
# document -- the new content belongs right after it).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*fixed by removing it*") {
        $targetIndex = $i
    }
}
if ($targetIndex -eq 0) {
    $targetIndex = $d.Paragraphs.Count
}

$r = $d.Paragraphs.Item($targetIndex).Range
$r.Collapse(0)

# Insert five blank paragraphs right after the target paragraph; we will
# stamp each one with its final OOXML content below.
for ($i = 0; $i -lt 5; $i++) {
    $r.InsertParagraphAfter() | Out-Null
    $r.Collapse(0)
    $r.Move(1, 1) | Out-Null
}

$pkgOpen = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$p1 = '<w:p/>'

$p2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>17/08/2020, 10:26:21 - Template contains errors.: Template parameters must be a map</w:t></w:r></w:p>'

$p3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Template missing required/referenced parameter definition in parameter section ' + [char]0x2013 + ' fixed by adding missing parameters</w:t></w:r></w:p>'

$p4 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'

$p5 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr></w:p>'

$newParas = @($p1, $p2, $p3, $p4, $p5)

for ($i = 0; $i -lt 5; $i++) {
    $p = $d.Paragraphs.Item($targetIndex + 1 + $i)
    $pr = $p.Range
    $pr.InsertXML($pkgOpen + $newParas[$i] + $pkgClose)
}

$d.Save()
